# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock) sheet gains three trailing columns - date,
# legislator_name, legislator_id - populated with the same value on every
# data row (this filing's date, the legislator's name and id).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 5          # existing data rows are 2..5 (row 1 is the header)
$dateCol = 8           # H
$nameCol = 9           # I
$idCol   = 10          # J

# --- header row -----------------------------------------------------------
$ws.Cells.Item(1, $dateCol).Value = "date"
$ws.Cells.Item(1, $nameCol).Value = "legislator_name"
$ws.Cells.Item(1, $idCol).Value   = "legislator_id"

# match the bold/bordered header formatting already used by columns B:G
$ws.Range("G1").Copy()
$ws.Range($ws.Cells.Item(1, $dateCol), $ws.Cells.Item(1, $idCol)).PasteSpecial(-4122)

# --- data rows --------------------------------------------------------------
# Pre-format the date column as text so the "2011-11-23" literal is kept as
# a plain string instead of being auto-converted into a date serial value.
$ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol)).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dateCol).Value = "2011-11-23"
    $ws.Cells.Item($r, $nameCol).Value = "翁重鈞"
    $ws.Cells.Item($r, $idCol).Value   = 551
}
